$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns to reflect the latest scrape.
# D-column values are forced to text (NumberFormat "@") so strings that look
# like numbers (e.g. "495.64", "1.00") are not auto-converted to numeric cells,
# then the style is reset to "Normal" so no residual text-format style sticks.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.329.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.33%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.291.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.88%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "495.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.83%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.528"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.291.70"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.95%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0943"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.36%  "

$ws.Range("E11").Value = "  +0.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.321"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.05%  "

$ws.Range("E13").Value = "  -2.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.696.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "54.265.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.41%  "

$ws.Range("E17").Value = "  -1.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.280.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "303.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.93%  "

$ws.Range("E25").Value = "  +0.40%  "

$ws.Range("E26").Value = "  +1.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.394.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.30%  "

$ws.Range("E28").Value = "  +2.21%  "

$ws.Range("E29").Value = "  +1.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "165.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.24%  "

$ws.Range("E31").Value = "  -1.70%  "

$ws.Range("E32").Value = "  -2.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.83%  "

$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("E36").Value = "  +1.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.22%  "

$ws.Range("E38").Value = "  +2.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.876"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.375"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "125.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0890"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.546"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "237.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.62%  "

$ws.Range("E50").Value = "  +1.47%  "

$ws.Range("E51").Value = "  +0.30%  "
